$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a flat "history" table (one match per row), sorted by date.
# Two new finished matches need to be inserted right before the existing
# row 262 (they happened slightly earlier than the match currently on row
# 262), so every existing row from 262 to 268 shifts down by two rows
# (-> 264..270), and the two vacated rows (262, 263) are filled with the
# newly reported matches. The closing odds (M..U) of the matches that used
# to sit on rows 262-264 are also refreshed to their final settled values.
# ---------------------------------------------------------------------------

# xlPasteFormats = -4122, xlPasteAll = -4104

# Shift existing rows 262..268 down to 264..270 (process bottom-up so the
# source of each copy is untouched until after it has been read).
for ($r = 268; $r -ge 262; $r--) {
    $d = $r + 2

    # Column A (id, bold/bordered style + number) - copy format, then value.
    $ws.Range("A$r").Copy()
    $ws.Range("A$d").PasteSpecial(-4122)
    $ws.Range("A$d").Value2 = $ws.Range("A$r").Value2

    # Columns B:C (text) - safe to paste whole (no date format involved).
    $ws.Range("B$r`:C$r").Copy()
    $ws.Range("B$d`:C$d").PasteSpecial(-4104)

    # Column D (date/time) - copy format, then value (avoids spurious
    # alternate-locale number-format being created).
    $ws.Range("D$r").Copy()
    $ws.Range("D$d").PasteSpecial(-4122)
    $ws.Range("D$d").Value2 = $ws.Range("D$r").Value2

    # Columns E:F (text).
    $ws.Range("E$r`:F$r").Copy()
    $ws.Range("E$d`:F$d").PasteSpecial(-4104)

    # Columns G:AB (plain numbers).
    $ws.Range("G$d`:AB$d").Value2 = $ws.Range("G$r`:AB$r").Value2
}

# Refresh closing odds (columns M:U) for the three matches that are now on
# rows 264-266 (previously 262-264) with their settled values.
$ws.Range("M264").Value2 = 1.533
$ws.Range("N264").Value2 = 4.333
$ws.Range("O264").Value2 = 5.5
$ws.Range("P264").Value2 = -1
$ws.Range("Q264").Value2 = 1.9
$ws.Range("R264").Value2 = 1.95
$ws.Range("S264").Value2 = 3
$ws.Range("T264").Value2 = 2
$ws.Range("U264").Value2 = 1.85

$ws.Range("M265").Value2 = 1.95
$ws.Range("N265").Value2 = 3.3
$ws.Range("O265").Value2 = 3.8
$ws.Range("P265").Value2 = -0.5
$ws.Range("Q265").Value2 = 1.975
$ws.Range("R265").Value2 = 1.875
$ws.Range("S265").Value2 = 2.25
$ws.Range("T265").Value2 = 1.95
$ws.Range("U265").Value2 = 1.9

$ws.Range("M266").Value2 = 3.6
$ws.Range("N266").Value2 = 3.5
$ws.Range("O266").Value2 = 2
$ws.Range("P266").Value2 = 0.5
$ws.Range("Q266").Value2 = 1.825
$ws.Range("R266").Value2 = 2.025
$ws.Range("S266").Value2 = 2.5
$ws.Range("T266").Value2 = 1.95
$ws.Range("U266").Value2 = 1.9

# ---------------------------------------------------------------------------
# New row 262: Piast Gliwice 2-0 Warta Poznan
# ---------------------------------------------------------------------------
$ws.Range("A262").Value2 = 260
$ws.Range("B262").Value2 = 6876422
$ws.Range("C262").Value = "Poland Ekstraklasa"
$ws.Range("D262").Value2 = 45408.54166666666
$ws.Range("E262").Value = "Piast Gliwice"
$ws.Range("F262").Value = "Warta Poznan"
$ws.Range("G262").Value2 = 2
$ws.Range("H262").Value2 = 0
$ws.Range("I262").Value = "H"
$ws.Range("J262").Value2 = 1.909
$ws.Range("K262").Value2 = 3.1
$ws.Range("L262").Value2 = 4.5
$ws.Range("M262").Value2 = 1.833
$ws.Range("N262").Value2 = 2.875
$ws.Range("O262").Value2 = 5.25
$ws.Range("P262").Value2 = -0.5
$ws.Range("Q262").Value2 = 1.875
$ws.Range("R262").Value2 = 1.975
$ws.Range("S262").Value2 = 1.75
$ws.Range("T262").Value2 = 1.8
$ws.Range("U262").Value2 = 2.05
$ws.Range("V262").Value2 = 0.833
$ws.Range("W262").Value2 = -1
$ws.Range("X262").Value2 = -1
$ws.Range("Y262").Value2 = 0.875
$ws.Range("Z262").Value2 = -1
$ws.Range("AA262").Value2 = 0.4
$ws.Range("AB262").Value2 = -0.5

# ---------------------------------------------------------------------------
# New row 263: Jagiellonia Bialystok 2-2 Pogon Szczecin
# ---------------------------------------------------------------------------
$ws.Range("A263").Value2 = 261
$ws.Range("B263").Value2 = 6921816
$ws.Range("C263").Value = "Poland Ekstraklasa"
$ws.Range("D263").Value2 = 45408.64583333334
$ws.Range("E263").Value = "Jagiellonia Bialystok"
$ws.Range("F263").Value = "Pogon Szczecin"
$ws.Range("G263").Value2 = 2
$ws.Range("H263").Value2 = 2
$ws.Range("I263").Value = "D"
$ws.Range("J263").Value2 = 2.1
$ws.Range("K263").Value2 = 3.75
$ws.Range("L263").Value2 = 3.1
$ws.Range("M263").Value2 = 1.909
$ws.Range("N263").Value2 = 4
$ws.Range("O263").Value2 = 3.4
$ws.Range("P263").Value2 = -0.5
$ws.Range("Q263").Value2 = 1.925
$ws.Range("R263").Value2 = 1.925
$ws.Range("S263").Value2 = 3.25
$ws.Range("T263").Value2 = 1.925
$ws.Range("U263").Value2 = 1.925
$ws.Range("V263").Value2 = -1
$ws.Range("W263").Value2 = 3
$ws.Range("X263").Value2 = -1
$ws.Range("Y263").Value2 = -1
$ws.Range("Z263").Value2 = 0.925
$ws.Range("AA263").Value2 = 0.925
$ws.Range("AB263").Value2 = -1

# Apply the same id/date cell style as the rest of the table to the two new
# rows (copy format only from row 264, which already carries the original
# style after the shift above).
$ws.Range("A264").Copy()
$ws.Range("A262:A263").PasteSpecial(-4122)
$ws.Range("D264").Copy()
$ws.Range("D262:D263").PasteSpecial(-4122)

Write-Host "Done"
